$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F2").Value = 0.8114832535885167
$ws.Range("G2").Value = 0.5619615639496355
$ws.Range("H2").Value = 0.6640563821456538
$ws.Range("E3").Value = 0.8426984916293718
$ws.Range("F3").Value = 0.7143950995405819
$ws.Range("G3").Value = 0.6182902584493042
$ws.Range("H3").Value = 0.6628774422735346
$ws.Range("E4").Value = 0.8417039615448367
$ws.Range("F4").Value = 0.7076461769115442
$ws.Range("G4").Value = 0.6255798542080848
$ws.Range("H4").Value = 0.6640872317973971
$ws.Range("E5").Value = 0.8417039615448367
$ws.Range("F5").Value = 0.7076461769115442
$ws.Range("G5").Value = 0.6255798542080848
$ws.Range("H5").Value = 0.6640872317973971
$ws.Range("E6").Value = 0.8417039615448367
$ws.Range("F6").Value = 0.7076461769115442
$ws.Range("G6").Value = 0.6255798542080848
$ws.Range("H6").Value = 0.6640872317973971
$ws.Range("F7").Value = 0.8117195004803074
$ws.Range("G7").Value = 0.559973492379059
$ws.Range("H7").Value = 0.6627450980392158
$ws.Range("E8").Value = 0.8410409414884801
$ws.Range("F8").Value = 0.7096036585365854
$ws.Range("G8").Value = 0.6169648774022531
$ws.Range("H8").Value = 0.6600496277915632
$ws.Range("F9").Value = 0.703288490284006
$ws.Range("G9").Value = 0.6235917826375083
$ws.Range("H9").Value = 0.6610467158412363
$ws.Range("F10").Value = 0.703288490284006
$ws.Range("G10").Value = 0.6235917826375083
$ws.Range("H10").Value = 0.6610467158412363
$ws.Range("F11").Value = 0.703288490284006
$ws.Range("G11").Value = 0.6235917826375083
$ws.Range("H11").Value = 0.6610467158412363
$ws.Range("F12").Value = 0.8155619596541787
$ws.Range("G12").Value = 0.562624254473161
$ws.Range("H12").Value = 0.6658823529411764
$ws.Range("F13").Value = 0.7075399847677075
$ws.Range("G13").Value = 0.6156394963552021
$ws.Range("H13").Value = 0.6583982990786675
$ws.Range("F14").Value = 0.7001499250374813
$ws.Range("G14").Value = 0.6189529489728297
$ws.Range("H14").Value = 0.657052409426662
$ws.Range("F15").Value = 0.7001499250374813
$ws.Range("G15").Value = 0.6189529489728297
$ws.Range("H15").Value = 0.657052409426662
$ws.Range("F16").Value = 0.7001499250374813
$ws.Range("G16").Value = 0.6189529489728297
$ws.Range("H16").Value = 0.657052409426662
$ws.Range("F17").Value = 0.8155619596541787
$ws.Range("G17").Value = 0.562624254473161
$ws.Range("H17").Value = 0.6658823529411764
$ws.Range("F18").Value = 0.7083015993907082
$ws.Range("G18").Value = 0.6163021868787276
$ws.Range("H18").Value = 0.659107016300496
$ws.Range("F19").Value = 0.7012012012012012
$ws.Range("G19").Value = 0.6189529489728297
$ws.Range("H19").Value = 0.6575149595212954
$ws.Range("F20").Value = 0.7012012012012012
$ws.Range("G20").Value = 0.6189529489728297
$ws.Range("H20").Value = 0.6575149595212954
$ws.Range("F21").Value = 0.7012012012012012
$ws.Range("G21").Value = 0.6189529489728297
$ws.Range("H21").Value = 0.6575149595212954
$ws.Range("F22").Value = 0.8159922928709056
$ws.Range("G22").Value = 0.56129887342611
$ws.Range("H22").Value = 0.6650961915979583
$ws.Range("F23").Value = 0.708649468892261
$ws.Range("G23").Value = 0.6189529489728297
$ws.Range("H23").Value = 0.6607711354793067
$ws.Range("F24").Value = 0.7030075187969925
$ws.Range("G24").Value = 0.6196156394963552
$ws.Range("H24").Value = 0.6586826347305389
$ws.Range("F25").Value = 0.7030075187969925
$ws.Range("G25").Value = 0.6196156394963552
$ws.Range("H25").Value = 0.6586826347305389
$ws.Range("F26").Value = 0.7030075187969925
$ws.Range("G26").Value = 0.6196156394963552
$ws.Range("H26").Value = 0.6586826347305389
